# Add the author's affiliation paragraph right after the "Edison Achalma"
# paragraph in the title/author block (the one styled "Author").
#
# There are several "Edison Achalma" occurrences in the document, so we
# first locate the specific paragraph (style "Author", text exactly
# "Edison Achalma") and then scope the Find/Replace to that paragraph's
# own Range so only that occurrence is touched. Using "^p" in the
# replacement text inserts a real paragraph break, and the new paragraph
# inherits the "Author" style from the paragraph it was split from.

$d = $word.ActiveDocument

$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Style.NameLocal -eq "Author" -and $p.Range.Text.Trim() -eq "Edison Achalma") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $r = $target.Range
    $r.Find.Execute(
        "Edison Achalma",
        $true,
        $false,
        $false,
        $false,
        $false,
        $true,
        1,
        $false,
        "Edison Achalma^pEscuela Profesional de Economía, Universidad Nacional de San Cristóbal de Huamanga",
        1
    )
}
